# Updates cryptos list values/links per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text even for numeric-looking values,
# matching the original workbook formatting (inline text strings, no numeric cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.543.14"
$ws.Range("E2").Value = "  +12.57%  "
$ws.Range("D3").Value = "1.838.65"
$ws.Range("E3").Value = "  +9.63%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "231.77"
$ws.Range("E5").Value = "  +5.31%  "
$ws.Range("D6").Value = "0.577"
$ws.Range("E6").Value = "  +8.82%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +7.09%  "
$ws.Range("D9").Value = "47.07"
$ws.Range("E9").Value = "  +6.24%  "
$ws.Range("D10").Value = "0.289"
$ws.Range("E10").Value = "  +9.04%  "
$ws.Range("D11").Value = "0.0682"
$ws.Range("E11").Value = "  +4.80%  "
$ws.Range("D12").Value = "0.0936"
$ws.Range("E12").Value = "  +3.41%  "
$ws.Range("D13").Value = "2.103.31"
$ws.Range("E13").Value = "  +9.67%  "
$ws.Range("D14").Value = "1.839.23"
$ws.Range("E14").Value = "  +9.81%  "
$ws.Range("E15").Value = "  +6.76%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "34.483.24"
$ws.Range("E16").Value = "  +12.34%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "10.38"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("D18").Value = "4.31"
$ws.Range("E18").Value = "  +7.02%  "
$ws.Range("D19").Value = "70.48"
$ws.Range("E19").Value = "  +5.93%  "
$ws.Range("D20").Value = "261.08"
$ws.Range("E20").Value = "  +7.29%  "
$ws.Range("D21").Value = "0.0₃0759"
$ws.Range("E21").Value = "  +4.54%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").Value = "  +6.11%  "
$ws.Range("D24").Value = "4.42"
$ws.Range("E24").Value = "  +3.76%  "
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").Value = "158.56"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "16.82"
$ws.Range("E27").Value = "  +5.85%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "0.119"
$ws.Range("E28").Value = "  +5.67%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "7.21"
$ws.Range("E29").Value = "  +7.56%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "3.89"
$ws.Range("E31").Value = "  +12.28%  "
$ws.Range("D32").Value = "0.0531"
$ws.Range("E32").Value = "  +7.05%  "
$ws.Range("E33").Value = "  +6.39%  "
$ws.Range("D34").Value = "3.61"
$ws.Range("E34").Value = "  +9.69%  "
$ws.Range("D35").Value = "1.561.53"
$ws.Range("E35").Value = "  +4.47%  "
$ws.Range("D36").Value = "1.82"
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("E37").Value = "  +5.81%  "
$ws.Range("D38").Value = "0.642"
$ws.Range("E38").Value = "  +7.59%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0192"
$ws.Range("E39").Value = "  +7.26%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "86.24"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.83"
$ws.Range("E41").Value = "  +5.77%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.926"
$ws.Range("E42").Value = "  +10.21%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.35"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "2.15"
$ws.Range("E44").Value = "  +8.16%  "
$ws.Range("B45").Value = "MinaProtocolToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D45").Value = "1.16"
$ws.Range("E45").Value = "  +179.27%  "
$ws.Range("D46").Value = "0.0527"
$ws.Range("E46").Value = "  +5.40%  "
$ws.Range("E47").Value = "  +6.40%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.996.29"
$ws.Range("E48").Value = "  +10.15%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "12.52"
$ws.Range("E49").Value = "  +22.68%  "
$ws.Range("D50").Value = "5.93"
$ws.Range("E50").Value = "  +6.94%  "
$ws.Range("E51").Value = "  +0.05%  "

# Remove the temporary text-number-format so cell styling matches the original (default) style.
$ws.Range("D2:D51").ClearFormats()

